# Auto-generated edit script: updates cryptos price/volume table
# to match the commit's scraped snapshot (Fri Aug 23 03:53:37 UTC 2024).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "60.501.28"
$ws.Range("E2").Value = "  +0.48%  "

# Row 3
$ws.Range("D3").Value = "2.634.48"
$ws.Range("E3").Value = "  +1.72%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "'581.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.39%  "

# Row 6
$ws.Range("D6").Value = "'144.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.20%  "

# Row 8
$ws.Range("D8").Value = "'0.598"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.17%  "

# Row 9
$ws.Range("E9").Value = "  -0.88%  "

# Row 10
$ws.Range("E10").Value = "  +0.80%  "

# Row 11
$ws.Range("E11").Value = "  +1.78%  "

# Row 12
$ws.Range("E12").Value = "  +3.46%  "

# Row 13
$ws.Range("D13").Value = "3.087.19"
$ws.Range("E13").Value = "  +1.16%  "

# Row 14
$ws.Range("D14").Value = "'26.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +14.15%  "

# Row 15
$ws.Range("D15").Value = "60.485.20"
$ws.Range("E15").Value = "  +0.46%  "

# Row 16
$ws.Range("E16").Value = "  +1.63%  "

# Row 17
$ws.Range("D17").Value = "2.628.17"
$ws.Range("E17").Value = "  +1.14%  "

# Row 18
$ws.Range("D18").Value = "'11.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.87%  "

# Row 19
$ws.Range("E19").Value = "  +1.74%  "

# Row 20
$ws.Range("D20").Value = "'346.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.38%  "

# Row 21
$ws.Range("E21").Value = "  -0.38%  "

# Row 22
$ws.Range("E22").Value = "  -0.07%  "

# Row 23
$ws.Range("E23").Value = "  -1.55%  "

# Row 24
$ws.Range("D24").Value = "'63.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.06%  "

# Row 25
$ws.Range("E25").Value = "  +0.03%  "

# Row 26
$ws.Range("E26").Value = "  +1.87%  "

# Row 27
$ws.Range("E27").Value = "  +6.18%  "

# Row 28
$ws.Range("D28").Value = "'2.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +14.22%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0796"
$ws.Range("E29").Value = "  +2.02%  "

# Row 30
$ws.Range("E30").Value = "  +5.27%  "

# Row 31
$ws.Range("D31").Value = "'169.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.73%  "

# Row 32
$ws.Range("E32").Value = "  +0.09%  "

# Row 33
$ws.Range("E33").Value = "  +0.94%  "

# Row 34
$ws.Range("D34").Value = "'1.08"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +12.62%  "

# Row 35
$ws.Range("D35").Value = "'4.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.48%  "

# Row 36
$ws.Range("D36").Value = "'1.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.77%  "

# Row 37
$ws.Range("E37").Value = "  +4.10%  "

# Row 38
$ws.Range("D38").Value = "'330.81"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +12.70%  "

# Row 39
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").Value = "'38.87"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.03%  "

# Row 40
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "'4.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.61%  "

# Row 41
$ws.Range("D41").Value = "'0.861"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.30%  "

# Row 42
$ws.Range("E42").Value = "  +6.66%  "

# Row 43
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "'20.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.10%  "

# Row 44
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'132.89"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.13%  "

# Row 45
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "'0.0996"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.00%  "

# Row 46
$ws.Range("D46").Value = "'20.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.31%  "

# Row 47
$ws.Range("D47").Value = "'1.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.35%  "

# Row 48
$ws.Range("D48").Value = "'0.0558"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.50%  "

# Row 49
$ws.Range("E49").Value = "  +0.75%  "

# Row 50
$ws.Range("E50").Value = "  +2.27%  "

# Row 51
$ws.Range("E51").Value = "  +0.77%  "
